$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Angptl1", "Tek", "ECs", 1, 0.3333333333333333, 0.3800063333333333, 1.140019, 0.002363300787678962, 0.002363300787678963, 3, 1, 64.221457, 192.664371, 0.8366610301096816, 0.8366610301096816, 24.40456039589433, 219.641043563049, 0.001977281671478503, 0.001977281671478503),
  @("ECs", "Angptl1", "Tek", "FAPs", 1, 0.3333333333333333, 0.3800063333333333, 1.140019, 0.002363300787678962, 0.002363300787678963, 3, 1, 10.582537, 31.747611, 0.1378666371208897, 0.1378666371208896, 4.021431082734333, 36.19287974460899, 0.0003258203321024482, 0.0003258203321024482),
  @("ECs", "Angptl1", "Tek", "sCs", 1, 0.3333333333333333, 0.3800063333333333, 1.140019, 0.002363300787678962, 0.002363300787678963, 3, 1, 1.955236666666667, 5.86571, 0.02547233276942866, 0.02547233276942866, 0.7430023164988888, 6.68702084849, [double]"6.01987840980114E-05", [double]"6.01987840980114E-05"),
  @("FAPs", "Angptl1", "Tek", "ECs", 3, 1, 159.049647, 477.148941, 0.9891470827288694, 0.9891470827288696, 3, 1, 64.221457, 192.664371, 0.8366610301096816, 0.8366610301096816, 10214.40006567568, 91929.60059108112, 0.8275808171659224, 0.8275808171659226),
  @("FAPs", "Angptl1", "Tek", "FAPs", 3, 1, 159.049647, 477.148941, 0.9891470827288694, 0.9891470827288696, 3, 1, 10.582537, 31.747611, 0.1378666371208897, 0.1378666371208896, 1683.148774214439, 15148.33896792995, 0.1363703819137677, 0.1363703819137677),
  @("FAPs", "Angptl1", "Tek", "sCs", 3, 1, 159.049647, 477.148941, 0.9891470827288694, 0.9891470827288696, 3, 1, 1.955236666666667, 5.86571, 0.02547233276942866, 0.02547233276942866, 310.97970163479, 2798.81731471311, 0.02519588364917934, 0.02519588364917934),
  @("M2", "Angptl1", "Tek", "ECs", 3, 1, 0.450387, 1.351161, 0.00280100582146534, 0.00280100582146534, 3, 1, 64.221457, 192.664371, 0.8366610301096816, 0.8366610301096816, 28.924509353859, 260.320584184731, 0.002343492415930406, 0.002343492415930406),
  @("M2", "Angptl1", "Tek", "FAPs", 3, 1, 0.450387, 1.351161, 0.00280100582146534, 0.00280100582146534, 3, 1, 10.582537, 31.747611, 0.1378666371208897, 0.1378666371208896, 4.766237091819001, 42.896133826371, 0.0003861652531614615, 0.0003861652531614615),
  @("M2", "Angptl1", "Tek", "sCs", 3, 1, 0.450387, 1.351161, 0.00280100582146534, 0.00280100582146534, 3, 1, 1.955236666666667, 5.86571, 0.02547233276942866, 0.02547233276942866, 0.8806131765900002, 7.92551858931, [double]"7.134815237347202E-05", [double]"7.134815237347202E-05"),
  @("sCs", "Angptl1", "Tek", "ECs", 3, 1, 0.9146986666666667, 2.744096, 0.005688610661986064, 0.005688610661986065, 3, 1, 64.221457, 192.664371, 0.8366610301096816, 0.8366610301096816, 58.74328108929067, 528.689529803616, 0.004759438856350179, 0.00475943885635018),
  @("sCs", "Angptl1", "Tek", "FAPs", 3, 1, 0.9146986666666667, 2.744096, 0.005688610661986064, 0.005688610661986065, 3, 1, 10.582537, 31.747611, 0.1378666371208897, 0.1378666371208896, 9.679832483850667, 87.118492354656, 0.0007842696218580567, 0.0007842696218580567),
  @("sCs", "Angptl1", "Tek", "sCs", 3, 1, 0.9146986666666667, 2.744096, 0.005688610661986064, 0.005688610661986065, 3, 1, 1.955236666666667, 5.86571, 0.02547233276942866, 0.02547233276942866, 1.788452372017778, 16.09607134816, 0.0001449021837778289, 0.0001449021837778289)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 2
  $row = $data[$i]
  for ($j = 0; $j -lt $row.Length; $j++) {
    $c = $j + 1
    $ws.Cells.Item($r, $c).Value = $row[$j]
  }
}

